$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
# row 19
$ws.Cells.Item(19, 8).Value = 3482.875  # H19: was 4499.4
$ws.Cells.Item(19, 9).Value = 1343.5  # I19: was 1997.5
$ws.Cells.Item(19, 10).Value = 5622.25  # J19: was 6167.3335
$ws.Cells.Item(19, 11).Value = 1343.5  # K19: was 1997.5
$ws.Cells.Item(19, 12).Value = 5622.25  # L19: was 6167.3335
$ws.Cells.Item(19, 13).Value = -1168.5  # M19: was -1822.5
$ws.Cells.Item(19, 14).Value = -5972.25  # N19: was -6517.3335
# row 86
$ws.Cells.Item(86, 8).Value = 2488.7778  # H86: was 2687
$ws.Cells.Item(86, 9).Value = 983.1667  # I86: was 999.2
$ws.Cells.Item(86, 11).Value = 983.1667  # K86: was 999.2
$ws.Cells.Item(86, 13).Value = 139.8333  # M86: was 123.8
# row 89
$ws.Cells.Item(89, 8).Value = 2488.7778  # H89: was 2687
$ws.Cells.Item(89, 9).Value = 983.1667  # I89: was 999.2
$ws.Cells.Item(89, 11).Value = 4915.8335  # K89: was 4996
$ws.Cells.Item(89, 13).Value = 700.1665000000003  # M89: was 620
# row 97
$ws.Cells.Item(97, 8).Value = 1200  # H97: was 0
$ws.Cells.Item(97, 10).Value = 1200  # J97: was 0
$ws.Cells.Item(97, 12).Value = 3600  # L97: was 0
$ws.Cells.Item(97, 14).Value = -4592  # N97: was None
# row 112
$ws.Cells.Item(112, 8).Value = 2315.2354  # H112: was 4408.8335
$ws.Cells.Item(112, 10).Value = 2334.9375  # J112: was 4550.5293
$ws.Cells.Item(112, 12).Value = 7004.8125  # L112: was 13651.5879
$ws.Cells.Item(112, 14).Value = -9220.8125  # N112: was -15867.5879
# row 113
$ws.Cells.Item(113, 8).Value = 13882.294  # H113: was 14637.5
$ws.Cells.Item(113, 9).Value = 15399.934  # I113: was 16371.429
$ws.Cells.Item(113, 11).Value = 15399.934  # K113: was 16371.429
$ws.Cells.Item(113, 13).Value = -12145.934  # M113: was -13117.429
# row 132
$ws.Cells.Item(132, 8).Value = 1011.9286  # H132: was 1082.8718
$ws.Cells.Item(132, 9).Value = 917.525  # I132: was 984.6486
$ws.Cells.Item(132, 11).Value = 2752.575  # K132: was 2953.9458
$ws.Cells.Item(132, 13).Value = -222.5749999999998  # M132: was -423.9458
# row 137
$ws.Cells.Item(137, 8).Value = 1742.24  # H137: was 1884.238
$ws.Cells.Item(137, 9).Value = 1498.5385  # I137: was 1649.2
$ws.Cells.Item(137, 10).Value = 2006.25  # J137: was 2097.9092
$ws.Cells.Item(137, 11).Value = 4495.6155  # K137: was 4947.6
$ws.Cells.Item(137, 12).Value = 6018.75  # L137: was 6293.7276
$ws.Cells.Item(137, 13).Value = -1945.6155  # M137: was -2397.6
$ws.Cells.Item(137, 14).Value = -11118.75  # N137: was -11393.7276
# row 138
$ws.Cells.Item(138, 8).Value = 2509.25  # H138: was 2535.3518
$ws.Cells.Item(138, 9).Value = 3239.4211  # I138: was 3249.1052
$ws.Cells.Item(138, 10).Value = 2134.2974  # J138: was 2147.8857
$ws.Cells.Item(138, 11).Value = 9718.263300000001  # K138: was 9747.3156
$ws.Cells.Item(138, 12).Value = 6402.8922  # L138: was 6443.657099999999
$ws.Cells.Item(138, 13).Value = -4578.263300000001  # M138: was -4607.3156
$ws.Cells.Item(138, 14).Value = -16682.8922  # N138: was -16723.6571
# row 141
$ws.Cells.Item(141, 8).Value = 3114.8845  # H141: was 3123.3076
$ws.Cells.Item(141, 9).Value = 2408.5454  # I141: was 2418.5
$ws.Cells.Item(141, 11).Value = 7225.6362  # K141: was 7255.5
$ws.Cells.Item(141, 13).Value = -2045.6362  # M141: was -2075.5

# Sheet 2
$ws = $wb.Worksheets.Item(2)
# row 32
$ws.Cells.Item(32, 8).Value = 2489.8472  # H32: was 2576.9023
$ws.Cells.Item(32, 9).Value = 1699.8029  # I32: was 1769.9265
$ws.Cells.Item(32, 11).Value = 1699.8029  # K32: was 1769.9265
$ws.Cells.Item(32, 13).Value = -1412.8029  # M32: was -1482.9265
# row 61
$ws.Cells.Item(61, 8).Value = 2562.2173  # H61: was 2178.4644
$ws.Cells.Item(61, 9).Value = 1360  # I61: was 1118.3
$ws.Cells.Item(61, 10).Value = 4816.375  # J61: was 4828.875
$ws.Cells.Item(61, 11).Value = 1360  # K61: was 1118.3
$ws.Cells.Item(61, 12).Value = 4816.375  # L61: was 4828.875
$ws.Cells.Item(61, 13).Value = -1148  # M61: was -906.3
$ws.Cells.Item(61, 14).Value = -5240.375  # N61: was -5252.875
# row 110
$ws.Cells.Item(110, 8).Value = 1784.4166  # H110: was 1658
$ws.Cells.Item(110, 9).Value = 1036.3636  # I110: was 1015.3077
$ws.Cells.Item(110, 11).Value = 1036.3636  # K110: was 1015.3077
$ws.Cells.Item(110, 13).Value = 1008.6364  # M110: was 1029.6923
# row 123
$ws.Cells.Item(123, 8).Value = 57000  # H123: was 0
$ws.Cells.Item(123, 10).Value = 57000  # J123: was 0
$ws.Cells.Item(123, 12).Value = 57000  # L123: was 0
$ws.Cells.Item(123, 14).Value = -66800  # N123: was None
# row 132
$ws.Cells.Item(132, 8).Value = 1448.2245  # H132: was 1347.3392
$ws.Cells.Item(132, 9).Value = 1077.3611  # I132: was 1000.2955
$ws.Cells.Item(132, 10).Value = 2475.2307  # J132: was 2619.8333
$ws.Cells.Item(132, 11).Value = 3232.0833  # K132: was 3000.8865
$ws.Cells.Item(132, 12).Value = 7425.6921  # L132: was 7859.499899999999
$ws.Cells.Item(132, 13).Value = -702.0833000000002  # M132: was -470.8864999999996
$ws.Cells.Item(132, 14).Value = -12485.6921  # N132: was -12919.4999
# row 136
$ws.Cells.Item(136, 8).Value = 2562.2173  # H136: was 2178.4644
$ws.Cells.Item(136, 9).Value = 1360  # I136: was 1118.3
$ws.Cells.Item(136, 10).Value = 4816.375  # J136: was 4828.875
$ws.Cells.Item(136, 11).Value = 4080  # K136: was 3354.9
$ws.Cells.Item(136, 12).Value = 14449.125  # L136: was 14486.625
$ws.Cells.Item(136, 13).Value = -1530  # M136: was -804.8999999999996
$ws.Cells.Item(136, 14).Value = -19549.125  # N136: was -19586.625

# Sheet 3
$ws = $wb.Worksheets.Item(3)
# row 20
$ws.Cells.Item(20, 8).Value = 1525.3529  # H20: was 1540.0588
$ws.Cells.Item(20, 9).Value = 1676.1538  # I20: was 1684.6154
$ws.Cells.Item(20, 10).Value = 1035.25  # J20: was 1070.25
$ws.Cells.Item(20, 11).Value = 1676.1538  # K20: was 1684.6154
$ws.Cells.Item(20, 12).Value = 1035.25  # L20: was 1070.25
$ws.Cells.Item(20, 13).Value = -1429.1538  # M20: was -1437.6154
$ws.Cells.Item(20, 14).Value = -1529.25  # N20: was -1564.25
# row 82
$ws.Cells.Item(82, 8).Value = 19331  # H82: was 16397.4
$ws.Cells.Item(82, 9).Value = 3996.5  # I82: was 7996.75
$ws.Cells.Item(82, 11).Value = 3996.5  # K82: was 7996.75
$ws.Cells.Item(82, 13).Value = -3613.5  # M82: was -7613.75
# row 85
$ws.Cells.Item(85, 8).Value = 19331  # H85: was 16397.4
$ws.Cells.Item(85, 9).Value = 3996.5  # I85: was 7996.75
$ws.Cells.Item(85, 11).Value = 3996.5  # K85: was 7996.75
$ws.Cells.Item(85, 13).Value = -2670.5  # M85: was -6670.75
# row 134
$ws.Cells.Item(134, 8).Value = 2002.7576  # H134: was 5766.793
$ws.Cells.Item(134, 9).Value = 2006.3667  # I134: was 6555.7085
$ws.Cells.Item(134, 10).Value = 1966.6666  # J134: was 1980
$ws.Cells.Item(134, 11).Value = 6019.1001  # K134: was 19667.1255
$ws.Cells.Item(134, 12).Value = 5899.9998  # L134: was 5940
$ws.Cells.Item(134, 13).Value = -3484.1001  # M134: was -17132.1255
$ws.Cells.Item(134, 14).Value = -10969.9998  # N134: was -11010

# Sheet 4
$ws = $wb.Worksheets.Item(4)
# row 31
$ws.Cells.Item(31, 8).Value = 1725.826  # H31: was 1667.2727
$ws.Cells.Item(31, 10).Value = 2626.25  # J31: was 2570.8572
$ws.Cells.Item(31, 12).Value = 2626.25  # L31: was 2570.8572
$ws.Cells.Item(31, 14).Value = -3216.25  # N31: was -3160.8572
# row 34
$ws.Cells.Item(34, 8).Value = 1725.826  # H34: was 1667.2727
$ws.Cells.Item(34, 10).Value = 2626.25  # J34: was 2570.8572
$ws.Cells.Item(34, 12).Value = 2626.25  # L34: was 2570.8572
$ws.Cells.Item(34, 14).Value = -3030.25  # N34: was -2974.8572
# row 107
$ws.Cells.Item(107, 8).Value = 513.1429000000001  # H107: was 510.33334
$ws.Cells.Item(107, 9).Value = 422.4  # I107: was 427.1111
$ws.Cells.Item(107, 10).Value = 740  # J107: was 760
$ws.Cells.Item(107, 11).Value = 422.4  # K107: was 427.1111
$ws.Cells.Item(107, 12).Value = 740  # L107: was 760
$ws.Cells.Item(107, 13).Value = 1497.6  # M107: was 1492.8889
$ws.Cells.Item(107, 14).Value = -4580  # N107: was -4600

# Sheet 5
$ws = $wb.Worksheets.Item(5)
# row 11
$ws.Cells.Item(11, 8).Value = 450  # H11: was 582.5
$ws.Cells.Item(11, 9).Value = 450  # I11: was 474
$ws.Cells.Item(11, 10).Value = 0  # J11: was 799.5
$ws.Cells.Item(11, 11).Value = 1350  # K11: was 1422
$ws.Cells.Item(11, 12).Value = 0  # L11: was 2398.5
$ws.Cells.Item(11, 13).ClearContents()  # M11: was -1282
$ws.Cells.Item(11, 14).Value = -1210  # N11: was -2678.5
# row 33
$ws.Cells.Item(33, 8).Value = 86.2  # H33: was 92.59999999999999
$ws.Cells.Item(33, 10).Value = 94  # J33: was 102
$ws.Cells.Item(33, 12).Value = 564  # L33: was 612
$ws.Cells.Item(33, 14).Value = -1130  # N33: was -1178
# row 50
$ws.Cells.Item(50, 8).Value = 166783490  # H50: was 90973840
$ws.Cells.Item(50, 9).Value = 348884.5  # I50: was 233489.67
$ws.Cells.Item(50, 10).Value = 250000780  # J50: was 125001470
$ws.Cells.Item(50, 11).Value = 1046653.5  # K50: was 700469.01
$ws.Cells.Item(50, 12).Value = 750002340  # L50: was 375004410
$ws.Cells.Item(50, 13).Value = -1046172.5  # M50: was -699988.01
$ws.Cells.Item(50, 14).Value = -750003302  # N50: was -375005372
# row 52
$ws.Cells.Item(52, 8).Value = 1000  # H52: was 997.5
$ws.Cells.Item(52, 10).Value = 1000  # J52: was 997.5
$ws.Cells.Item(52, 12).Value = 3000  # L52: was 2992.5
$ws.Cells.Item(52, 14).Value = -3532  # N52: was -3524.5
# row 53
$ws.Cells.Item(53, 8).Value = 166783490  # H53: was 90973840
$ws.Cells.Item(53, 9).Value = 348884.5  # I53: was 233489.67
$ws.Cells.Item(53, 10).Value = 250000780  # J53: was 125001470
$ws.Cells.Item(53, 11).Value = 1046653.5  # K53: was 700469.01
$ws.Cells.Item(53, 12).Value = 750002340  # L53: was 375004410
$ws.Cells.Item(53, 13).Value = -1046172.5  # M53: was -699988.01
$ws.Cells.Item(53, 14).Value = -750003302  # N53: was -375005372
# row 75
$ws.Cells.Item(75, 8).Value = 1700  # H75: was 0
$ws.Cells.Item(75, 10).Value = 1700  # J75: was 0
$ws.Cells.Item(75, 12).Value = 5100  # L75: was 0
$ws.Cells.Item(75, 14).Value = -7096  # N75: was None
# row 78
$ws.Cells.Item(78, 8).Value = 1700  # H78: was 0
$ws.Cells.Item(78, 10).Value = 1700  # J78: was 0
$ws.Cells.Item(78, 12).Value = 15300  # L78: was 0
$ws.Cells.Item(78, 14).Value = -25284  # N78: was None
# row 87
$ws.Cells.Item(87, 8).Value = 13374  # H87: was 14001.125
$ws.Cells.Item(87, 9).Value = 6398.6  # I87: was 7402
$ws.Cells.Item(87, 11).Value = 19195.8  # K87: was 22206
$ws.Cells.Item(87, 13).Value = -17947.8  # M87: was -20958
# row 90
$ws.Cells.Item(90, 8).Value = 13374  # H90: was 14001.125
$ws.Cells.Item(90, 9).Value = 6398.6  # I90: was 7402
$ws.Cells.Item(90, 11).Value = 57587.4  # K90: was 66618
$ws.Cells.Item(90, 13).Value = -51347.4  # M90: was -60378
# row 103
$ws.Cells.Item(103, 8).Value = 2050.6428  # H103: was 1783.909
$ws.Cells.Item(103, 9).Value = 1546.5714  # I103: was 1324.8
$ws.Cells.Item(103, 10).Value = 2554.7144  # J103: was 2166.5
$ws.Cells.Item(103, 11).Value = 4639.7142  # K103: was 3974.4
$ws.Cells.Item(103, 12).Value = 7664.1432  # L103: was 6499.5
$ws.Cells.Item(103, 13).Value = -3760.7142  # M103: was -3095.4
$ws.Cells.Item(103, 14).Value = -9422.143199999999  # N103: was -8257.5
# row 108
$ws.Cells.Item(108, 8).Value = 2004.2  # H108: was 2007
$ws.Cells.Item(108, 9).Value = 2004.2  # I108: was 2007
$ws.Cells.Item(108, 11).Value = 6012.6  # K108: was 6021
$ws.Cells.Item(108, 13).Value = -3132.6  # M108: was -3141
# row 109
$ws.Cells.Item(109, 8).Value = 1366.4286  # H109: was 1429.8334
$ws.Cells.Item(109, 9).Value = 927.5  # I109: was 915.8
$ws.Cells.Item(109, 11).Value = 2782.5  # K109: was 2747.4
$ws.Cells.Item(109, 13).Value = -1742.5  # M109: was -1707.4
# row 114
$ws.Cells.Item(114, 8).Value = 3533  # H114: was 2391
$ws.Cells.Item(114, 9).Value = 0  # I114: was 678
$ws.Cells.Item(114, 11).Value = 0  # K114: was 2034
$ws.Cells.Item(114, 13).ClearContents()  # M114: was 1220
# row 131
$ws.Cells.Item(131, 8).Value = 1665.02  # H131: was 2800.3
$ws.Cells.Item(131, 10).Value = 1732.3937  # J131: was 2940.1382
$ws.Cells.Item(131, 12).Value = 5197.1811  # L131: was 8820.4146
$ws.Cells.Item(131, 14).Value = -15277.1811  # N131: was -18900.4146

# Sheet 6
$ws = $wb.Worksheets.Item(6)
# row 122
$ws.Cells.Item(122, 8).Value = 2016.5294  # H122: was 2036.3125
$ws.Cells.Item(122, 10).Value = 2358.7  # J122: was 2431.889
$ws.Cells.Item(122, 12).Value = 7076.099999999999  # L122: was 7295.667
$ws.Cells.Item(122, 14).Value = -11976.1  # N122: was -12195.667
# row 132
$ws.Cells.Item(132, 8).Value = 2602  # H132: was 2783.077
$ws.Cells.Item(132, 9).Value = 2359.7827  # I132: was 2537.8696
$ws.Cells.Item(132, 10).Value = 3994.75  # J132: was 4663
$ws.Cells.Item(132, 11).Value = 7079.348100000001  # K132: was 7613.6088
$ws.Cells.Item(132, 12).Value = 11984.25  # L132: was 13989
$ws.Cells.Item(132, 13).Value = -4549.348100000001  # M132: was -5083.6088
$ws.Cells.Item(132, 14).Value = -17044.25  # N132: was -19049
# row 135
$ws.Cells.Item(135, 8).Value = 29333.334  # H135: was 29000
$ws.Cells.Item(135, 10).Value = 29333.334  # J135: was 29000
$ws.Cells.Item(135, 12).Value = 29333.334  # L135: was 29000
$ws.Cells.Item(135, 14).Value = -39473.334  # N135: was -39140

# Sheet 7
$ws = $wb.Worksheets.Item(7)
# row 82
$ws.Cells.Item(82, 8).Value = 3568.5715  # H82: was 3663.3333
$ws.Cells.Item(82, 9).Value = 2750  # I82: was 2666.6667
$ws.Cells.Item(82, 11).Value = 2750  # K82: was 2666.6667
$ws.Cells.Item(82, 13).Value = -2389  # M82: was -2305.6667
# row 85
$ws.Cells.Item(85, 8).Value = 3568.5715  # H85: was 3663.3333
$ws.Cells.Item(85, 9).Value = 2750  # I85: was 2666.6667
$ws.Cells.Item(85, 11).Value = 2750  # K85: was 2666.6667
$ws.Cells.Item(85, 13).Value = -1502  # M85: was -1418.6667

# Sheet 8
$ws = $wb.Worksheets.Item(8)
# row 70
$ws.Cells.Item(70, 8).Value = 29400  # H70: was 29425
$ws.Cells.Item(70, 10).Value = 29400  # J70: was 29425
$ws.Cells.Item(70, 12).Value = 29400  # L70: was 29425
$ws.Cells.Item(70, 14).Value = -30030  # N70: was -30055
# row 73
$ws.Cells.Item(73, 8).Value = 29400  # H73: was 29425
$ws.Cells.Item(73, 10).Value = 29400  # J73: was 29425
$ws.Cells.Item(73, 12).Value = 29400  # L73: was 29425
$ws.Cells.Item(73, 14).Value = -31584  # N73: was -31609
# row 113
$ws.Cells.Item(113, 8).Value = 2299.6667  # H113: was 2750
$ws.Cells.Item(113, 10).Value = 2159.6  # J113: was 2666.6667
$ws.Cells.Item(113, 12).Value = 6478.799999999999  # L113: was 8000.000100000001
$ws.Cells.Item(113, 14).Value = -10818.8  # N113: was -12340.0001
# row 126
$ws.Cells.Item(126, 8).Value = 4835.2  # H126: was 5750.25
$ws.Cells.Item(126, 9).Value = 2621.7144  # I126: was 3200.4
$ws.Cells.Item(126, 11).Value = 7865.1432  # K126: was 9601.200000000001
$ws.Cells.Item(126, 13).Value = -5395.1432  # M126: was -7131.200000000001
# row 132
$ws.Cells.Item(132, 8).Value = 1392.0358  # H132: was 1527.9048
$ws.Cells.Item(132, 9).Value = 1179.4286  # I132: was 1224.1875
$ws.Cells.Item(132, 10).Value = 2029.8572  # J132: was 2499.8
$ws.Cells.Item(132, 11).Value = 3538.2858  # K132: was 3672.5625
$ws.Cells.Item(132, 12).Value = 6089.571599999999  # L132: was 7499.400000000001
$ws.Cells.Item(132, 13).Value = -1008.2858  # M132: was -1142.5625
$ws.Cells.Item(132, 14).Value = -11149.5716  # N132: was -12559.4
# row 136
$ws.Cells.Item(136, 8).Value = 2245.923  # H136: was 2329.88
$ws.Cells.Item(136, 9).Value = 1919.2222  # I136: was 1930.4445
$ws.Cells.Item(136, 10).Value = 2981  # J136: was 3357
$ws.Cells.Item(136, 11).Value = 5757.6666  # K136: was 5791.333500000001
$ws.Cells.Item(136, 12).Value = 8943  # L136: was 10071
$ws.Cells.Item(136, 13).Value = -3207.6666  # M136: was -3241.333500000001
$ws.Cells.Item(136, 14).Value = -14043  # N136: was -15171
